$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.076312065124512
$ws.Range("B1").Value = 1.746430993080139
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.907104134559631
$ws.Range("E1").Value = 1.153485059738159
